{"js": "// Applies the \"Added detail about how the brute force was done.\" edit.\n//\n// Four related changes (see XML diff):\n//  1. \"Paul Eccleston\" -> \"Paul \" + proofErr(spellStart) + \"Eccleston\" + proofErr(spellEnd)\n//  2. \"CS 465 Program 1\" -> \"CS 465 \" + proofErr(gramStart) + \"Program\" + proofErr(gramEnd) + \" 1\"\n//  3. Procedures paragraph gets a spellStart/spellEnd around \"openssl\" and the brute-force\n//     sentence is rewritten/expanded with new wording; the \"_GoBack\" bookmark moves into\n//     this paragraph (between \"...then we estimated\" and \" how long it would take...\").\n//  4. The trailing \"_GoBack\" bookmark is removed from the Conclusion paragraph (it moved\n//     up into the Procedures paragraph).\n\nconst OOXML_NS_WRAPPER = (bodyFragment) => `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">` +\n  `<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">` +\n  `<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">` +\n  `<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>` +\n  `</Relationships></pkg:xmlData></pkg:part>` +\n  `<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">` +\n  `<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">` +\n  `<w:body>${bodyFragment}</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>`;\n\nconst body = context.document.body;\n\n// ---- 1. \"Paul Eccleston\" -------------------------------------------------\nconst paulResults = body.search(\"Paul Eccleston\", { matchCase: true });\npaulResults.load(\"text\");\nawait context.sync();\nif (paulResults.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for 'Paul Eccleston', found \" + paulResults.items.length);\n}\npaulResults.items[0].insertOoxml(\n  OOXML_NS_WRAPPER(\n    '<w:p>' +\n      '<w:r><w:t xml:space=\"preserve\">Paul </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:t>Eccleston</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n    '</w:p>'\n  ),\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// ---- 2. \"CS 465 Program 1\" ------------------------------------------------\nconst titleResults = body.search(\"CS 465 Program 1\", { matchCase: true });\ntitleResults.load(\"text\");\nawait context.sync();\nif (titleResults.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for 'CS 465 Program 1', found \" + titleResults.items.length);\n}\ntitleResults.items[0].insertOoxml(\n  OOXML_NS_WRAPPER(\n    '<w:p>' +\n      '<w:r><w:t xml:space=\"preserve\">CS 465 </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:t>Program</w:t></w:r>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> 1</w:t></w:r>' +\n    '</w:p>'\n  ),\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// ---- 3. Remove the old \"_GoBack\" bookmark from the Conclusion paragraph ---\n// (Done BEFORE re-adding it in the Procedures paragraph below, since\n// deleteBookmark(\"_GoBack\") removes the first-in-document-order bookmark\n// with that name.)\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// ---- 4. Procedures paragraph ----------------------------------------------\nconst procOld = \"The encryption and decryption were done using bash shell script and openssl. A program was created to run ten iterations of encryption using the three methods, followed by ten decryptions of each. The brute force was done allowing the bash script to try to decrypt the file over a period of time, then calculating out based on the number of decades it would take for completion.  \";\nconst procResults = body.search(procOld, { matchCase: true });\nprocResults.load(\"text\");\nawait context.sync();\nif (procResults.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the Procedures paragraph text, found \" + procResults.items.length);\n}\nprocResults.items[0].insertOoxml(\n  OOXML_NS_WRAPPER(\n    '<w:p>' +\n      '<w:r><w:t xml:space=\"preserve\">The encryption and decryption were done using bash shell script and </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:t>openssl</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\">. A program was created to run ten iterations of encryption using the three methods, followed by ten decryptions of each. The brute force was done allowing the bash script to try to decrypt </w:t></w:r>' +\n      '<w:r><w:t>a file 1000 times with a wrong password, then we estimated</w:t></w:r>' +\n      '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n      '<w:bookmarkEnd w:id=\"0\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> how long it would take by calculating how long it would take to try every possible solution.</w:t></w:r>' +\n      '<w:r><w:t xml:space=\"preserve\">  </w:t></w:r>' +\n    '</w:p>'\n  ),\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Applies the \"Added detail about how the brute force was done.\" edit.\n#\n# Four related changes (see XML diff):\n#  1. \"Paul Eccleston\" -> \"Paul \" + proofErr(spellStart) + \"Eccleston\" + proofErr(spellEnd)\n#  2. \"CS 465 Program 1\" -> \"CS 465 \" + proofErr(gramStart) + \"Program\" + proofErr(gramEnd) + \" 1\"\n#  3. The old \"_GoBack\" bookmark is removed from the end of the Conclusion paragraph\n#     (done before step 4 re-creates it, so bookmark ids stay simple/clean).\n#  4. Procedures paragraph gets a spellStart/spellEnd around \"openssl\" and the brute-force\n#     sentence is rewritten/expanded with new wording; a new \"_GoBack\" bookmark is placed\n#     inside this paragraph (between \"...then we estimated\" and \" how long it would take...\").\n\n# Replaces the run content found via Find.Execute(searchText) with the supplied raw\n# OOXML run markup ($runsXml, e.g. \"<w:r>...</w:r><w:proofErr .../>...\") while preserving\n# the containing paragraph's <w:pPr/> (formatting). Works by inserting the replacement\n# XML directly after the matched text (so it merges into the existing paragraph instead\n# of becoming a brand-new paragraph), then clearing out the original matched text that is\n# now immediately in front of it.\nfunction Replace-RunsInParagraph($doc, $searchText, $runsXml) {\n    $pkgXml = \"<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>\" +\n              \"<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>\" +\n              \"<pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>\" +\n              \"<w:body><w:p>$runsXml</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\"\n\n    $rngInsert = $doc.Content\n    $rngInsert.Find.Text = $searchText\n    $found = $rngInsert.Find.Execute()\n    if (-not $found) {\n        throw \"Replace-RunsInParagraph: could not find text: $searchText\"\n    }\n    $rngInsert.InsertXML($pkgXml)\n\n    $rngClear = $doc.Content\n    $rngClear.Find.Text = $searchText\n    $found2 = $rngClear.Find.Execute()\n    if (-not $found2) {\n        throw \"Replace-RunsInParagraph: could not find original text after insert: $searchText\"\n    }\n    $rngClear.Text = \"\"\n}\n\n$d = $word.ActiveDocument\n\n# ---- 1. \"Paul Eccleston\" ---------------------------------------------------\nReplace-RunsInParagraph $d \"Paul Eccleston\" (\n    \"<w:r><w:t xml:space='preserve'>Paul </w:t></w:r>\" +\n    \"<w:proofErr w:type='spellStart'/>\" +\n    \"<w:r><w:t>Eccleston</w:t></w:r>\" +\n    \"<w:proofErr w:type='spellEnd'/>\"\n)\n\n# ---- 2. \"CS 465 Program 1\" -------------------------------------------------\nReplace-RunsInParagraph $d \"CS 465 Program 1\" (\n    \"<w:r><w:t xml:space='preserve'>CS 465 </w:t></w:r>\" +\n    \"<w:proofErr w:type='gramStart'/>\" +\n    \"<w:r><w:t>Program</w:t></w:r>\" +\n    \"<w:proofErr w:type='gramEnd'/>\" +\n    \"<w:r><w:t xml:space='preserve'> 1</w:t></w:r>\"\n)\n\n# ---- 3. Remove the old \"_GoBack\" bookmark from the Conclusion paragraph ---\n# (Must happen BEFORE step 4 re-adds \"_GoBack\" in the Procedures paragraph.)\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# ---- 4. Procedures paragraph ----------------------------------------------\n$procOld = \"The encryption and decryption were done using bash shell script and openssl. A program was created to run ten iterations of encryption using the three methods, followed by ten decryptions of each. The brute force was done allowing the bash script to try to decrypt the file over a period of time, then calculating out based on the number of decades it would take for completion.  \"\n\n$procNewRuns = \"<w:r><w:t xml:space='preserve'>The encryption and decryption were done using bash shell script and </w:t></w:r>\" +\n    \"<w:proofErr w:type='spellStart'/>\" +\n    \"<w:r><w:t>openssl</w:t></w:r>\" +\n    \"<w:proofErr w:type='spellEnd'/>\" +\n    \"<w:r><w:t xml:space='preserve'>. A program was created to run ten iterations of encryption using the three methods, followed by ten decryptions of each. The brute force was done allowing the bash script to try to decrypt </w:t></w:r>\" +\n    \"<w:r><w:t>a file 1000 times with a wrong password, then we estimated</w:t></w:r>\" +\n    \"<w:bookmarkStart w:id='0' w:name='_GoBack'/>\" +\n    \"<w:bookmarkEnd w:id='0'/>\" +\n    \"<w:r><w:t xml:space='preserve'> how long it would take by calculating how long it would take to try every possible solution.</w:t></w:r>\" +\n    \"<w:r><w:t xml:space='preserve'>  </w:t></w:r>\"\n\nReplace-RunsInParagraph $d $procOld $procNewRuns\n"}
